# Weekly update: a new "Espárragos" price record for
# Terminal Hortofrutícola Agro Chillán was published, so it gets
# inserted at the top of the data block (row 18, right after the header
# block that ends at row 17) and every existing record from row 18
# down to row 50 shifts down by one row (to rows 19-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18. This pushes the previous rows
# 18..50 down to 19..51 (and copies row 18's formatting, e.g. the
# date style in column D, onto the new blank row) - same as Excel's
# normal "Insert Sheet Rows" behavior.
$ws.Rows.Item(18).Insert()

# Fill the newly inserted row 18 with the new week's record.
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C18").Value = 'Ñuble'
$ws.Range("D18").Value = 45195
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = 'Espárragos'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1400
$ws.Range("N18").Value = '$/kilo'
$ws.Range("O18").Value = 'Región de Ñuble'
$ws.Range("P18").Value = 1400
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 'Hortaliza'
